# Worked on the project
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Settings sheet: replace the Orchestrator queue settings (rows 2-3, plus the
# blank spacer row 4) with a set of Acme Test environment URLs, keeping the
# existing "logF_BusinessProcessName" row right after the header.
# ---------------------------------------------------------------------------
$settings = $wb.Worksheets.Item("Settings")

# Drop the old OrchestratorQueueName / OrchestratorQueueFolder rows (2 and 3)
# and the blank spacer row that used to follow them (old row 4). Deleting the
# same row index three times removes what were rows 2, 3 and 4.
$settings.Rows.Item(2).EntireRow.Delete()
$settings.Rows.Item(2).EntireRow.Delete()
$settings.Rows.Item(2).EntireRow.Delete()

# Row 2 is now logF_BusinessProcessName / Framework / description - it keeps
# its original wrapped/description styling automatically (style carries over
# with the row), nothing further to do here.

# New Acme environment rows.
$settings.Range("A3").Value = "AcmeHome"
$settings.Range("B3").Value = "https://acme-test.uipath.com/"

$settings.Range("A4").Value = "AcmeUrl"
$settings.Range("B4").Value = "https://acme-test.uipath.com/login"

$settings.Range("A5").Value = "AcmeWorkItems"
$settings.Range("B5").Value = "https://acme-test.uipath.com/work-items"

$settings.Range("A6").Value = "AcmeAllVendors"
$settings.Range("B6").Value = "https://acme-test.uipath.com/vendors/search-by-name"

$settings.Range("A7").Value = "AcmeVendorsSearch"
$settings.Range("B7").Value = "https://acme-test.uipath.com/vendors/search"

# ---------------------------------------------------------------------------
# Constants sheet: MaxRetryNumber moves from 0 to 2.
# ---------------------------------------------------------------------------
$constants = $wb.Worksheets.Item("Constants")
$constants.Range("B2").Value = 2
$constants.Range("B3").Select()

# ---------------------------------------------------------------------------
# Settings becomes the active sheet/tab (was Assets before).
# ---------------------------------------------------------------------------
$settings.Activate()
$settings.Range("A7").Select()
